# "aggiunti valore indice Gulpease"
# Fill in the Gulpease-index row (row 11) with its computed numeric
# sub-scores, replacing the placeholder text values that were there
# before, and move the active selection to C12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Indice gulpease") gets real numbers instead of the old
# text placeholders ("126.30", "51.54", "158.53").
$ws.Range("B11").Value = 51
$ws.Range("C11").Value = 54
$ws.Range("D11").Value = 60
$ws.Range("E11").Value = 55
$ws.Range("F11").Value = 67

# Match the integer right-aligned number format used by the rest of
# the sheet's data cells.
$ws.Range("B11:G11").NumberFormat = "0"

# The author's selection ended up on C12 after entering the data.
$ws.Range("C12").Select() | Out-Null
